# Saldo_guide.xlsx update - refreshed balance export (20240717 -> 20240718 run)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export timestamp
$ws.Name = "IClientBalance-20240718-093257-"

# All rows' "date" column (G) moves from 45490 (2024-07-17) to 45491 (2024-07-18)
$ws.Range("G2:G275").Value = 45491

# A handful of rows also had their balance figures (columns D/E/H) revised
# Row 52: E/H change from 2680.09 -> 2583.42
$ws.Range("E52").Value = 2583.42
$ws.Range("H52").Value = 2583.42

# Row 58: D changes from 0 -> 15545.4, H changes from 132.84 -> 15678.24 (E stays 132.84)
$ws.Range("D58").Value = 15545.4
$ws.Range("H58").Value = 15678.24

# Row 118: E/H change from 19019.490000000002 -> 5.01
$ws.Range("E118").Value = 5.01
$ws.Range("H118").Value = 5.01

# Row 255: E/H change from 10168.82 -> 16673.439999999999
$ws.Range("E255").Value = 16673.44
$ws.Range("H255").Value = 16673.44
